$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
$ws.Columns.Item(7).ColumnWidth = 8 - 5/6
$ws.Columns.Item(10).ColumnWidth = 8 - 5/6
$ws.Columns.Item(11).ColumnWidth = 8 - 5/6
$ws.Columns.Item(15).ColumnWidth = 8 - 5/6
$ws.Columns.Item(17).ColumnWidth = 8 - 5/6
$ws.Columns.Item(22).ColumnWidth = 8 - 5/6
$ws.Columns.Item(27).ColumnWidth = 8 - 5/6
$ws.Columns.Item(28).ColumnWidth = 8 - 5/6
$ws.Columns.Item(31).ColumnWidth = 5 - 5/6

# --- Update data rows 2-5 with new values ---
$row2 = New-Object "object[,]" 1,34
$row2[0,0] = 45046.50694444445
$row2[0,1] = 4.928
$row2[0,2] = 5.344
$row2[0,3] = 0
$row2[0,4] = 5.978
$row2[0,5] = 10.253
$row2[0,6] = 2.511
$row2[0,7] = 8.567
$row2[0,8] = 4.041
$row2[0,9] = 1.748
$row2[0,10] = 4.53
$row2[0,11] = 5.048
$row2[0,12] = 5.182
$row2[0,13] = 0.791
$row2[0,14] = 3.479
$row2[0,15] = 3.593
$row2[0,16] = 1.503
$row2[0,17] = 1.042
$row2[0,18] = 0.417
$row2[0,19] = 43.098
$row2[0,20] = 7.946
$row2[0,21] = 4.36
$row2[0,22] = 6.645
$row2[0,23] = 2.748
$row2[0,24] = 0.484
$row2[0,25] = 2.765
$row2[0,26] = 1.536
$row2[0,27] = 3.068
$row2[0,28] = 2.75
$row2[0,29] = 5.612
$row2[0,30] = 0
$row2[0,31] = 5.035
$row2[0,32] = 2.082
$row2[0,33] = 3.31
$ws.Range("A2:AH2").Value = $row2

$row3 = New-Object "object[,]" 1,34
$row3[0,0] = 45046.51388888889
$row3[0,1] = 21.464
$row3[0,2] = 16.77
$row3[0,3] = 0.523
$row3[0,4] = 44.487
$row3[0,5] = 39.212
$row3[0,6] = 16.124
$row3[0,7] = 58.375
$row3[0,8] = 24.993
$row3[0,9] = 11.342
$row3[0,10] = 17.784
$row3[0,11] = 18.934
$row3[0,12] = 19.97
$row3[0,13] = 5.203
$row3[0,14] = 16.584
$row3[0,15] = 23.077
$row3[0,16] = 12.988
$row3[0,17] = 0.766
$row3[0,18] = 0.853
$row3[0,19] = 244.068
$row3[0,20] = 45.574
$row3[0,21] = 15.767
$row3[0,22] = 31.68
$row3[0,23] = 16.26
$row3[0,24] = 2.182
$row3[0,25] = 28.786
$row3[0,26] = 12.867
$row3[0,27] = 12.113
$row3[0,28] = 14.06
$row3[0,29] = 20.164
$row3[0,30] = 0
$row3[0,31] = 51.891
$row3[0,32] = 8.934
$row3[0,33] = 18.777
$ws.Range("A3:AH3").Value = $row3

$row4 = New-Object "object[,]" 1,34
$row4[0,0] = 45046.52083333334
$row4[0,1] = 11.269
$row4[0,2] = 8.879
$row4[0,3] = 0.218
$row4[0,4] = 23.081
$row4[0,5] = 20.659
$row4[0,6] = 8.332
$row4[0,7] = 36.871
$row4[0,8] = 12.981
$row4[0,9] = 5.958
$row4[0,10] = 9.398
$row4[0,11] = 9.95
$row4[0,12] = 10.476
$row4[0,13] = 2.722
$row4[0,14] = 8.665
$row4[0,15] = 12.01
$row4[0,16] = 6.677
$row4[0,17] = 0.502
$row4[0,18] = 0.5
$row4[0,19] = 124.062
$row4[0,20] = 23.906
$row4[0,21] = 8.282
$row4[0,22] = 16.761
$row4[0,23] = 8.528
$row4[0,24] = 1.14
$row4[0,25] = 17.047
$row4[0,26] = 6.631
$row4[0,27] = 6.358
$row4[0,28] = 7.41
$row4[0,29] = 10.633
$row4[0,30] = 0
$row4[0,31] = 32.981
$row4[0,32] = 4.723
$row4[0,33] = 9.767
$ws.Range("A4:AH4").Value = $row4

$row5 = New-Object "object[,]" 1,34
$row5[0,0] = 45046.52777777778
$row5[0,1] = 6.89
$row5[0,2] = 5.49
$row5[0,3] = 0.1
$row5[0,4] = 13.93
$row5[0,5] = 12.68
$row5[0,6] = 5.01
$row5[0,7] = 23.39
$row5[0,8] = 7.85
$row5[0,9] = 3.64
$row5[0,10] = 5.78
$row5[0,11] = 6.1
$row5[0,12] = 6.4
$row5[0,13] = 1.66
$row5[0,14] = 5.27
$row5[0,15] = 7.26
$row5[0,16] = 3.99
$row5[0,17] = 0.37
$row5[0,18] = 0.34
$row5[0,19] = 72.62
$row5[0,20] = 14.52
$row5[0,21] = 5.07
$row5[0,22] = 10.28
$row5[0,23] = 5.21
$row5[0,24] = 0.69
$row5[0,25] = 10.54
$row5[0,26] = 3.97
$row5[0,27] = 3.89
$row5[0,28] = 4.54
$row5[0,29] = 6.53
$row5[0,30] = 0
$row5[0,31] = 20.77
$row5[0,32] = 2.91
$row5[0,33] = 5.91
$ws.Range("A5:AH5").Value = $row5

# --- Remove row 6 (only 4 data rows remain: 2-5) ---
$ws.Rows.Item(6).Delete()
